# "Finish the intro part of the weicc docs #2"
#
# 1) The deck's cached "datetimeFigureOut" field (master, every content
#    slide layout, and the notes master) is bumped from 4/17/22 -> 5/14/22.
# 2) Slide 2 fixes the typo "Font End" -> "Front End" in the pipeline
#    diagram.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$newDate = "5/14/22"

# --- Slide master's own Date placeholder ---------------------------------
# (HeadersFooters.DateAndTime.Text is a no-op at the master/layout level in
# this host, so go through the placeholder shape's TextRange instead.)
foreach ($shp in $master.Shapes) {
    if ($shp.Name -like "*Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

# --- Every slide layout under the master ---------------------------------
foreach ($layout in $master.CustomLayouts) {
    foreach ($shp in $layout.Shapes) {
        if ($shp.Name -like "*Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Notes master ----------------------------------------------------------
# (Here it's the opposite: the placeholder shape's TextRange assignment
# doesn't stick, but HeadersFooters.DateAndTime.Text does.)
$nm = $p.NotesMaster
$nm.HeadersFooters.DateAndTime.Text = $newDate

# --- Slide 2: "Font End" -> "Front End" -----------------------------------
$s2 = $p.Slides.Item(2)
foreach ($shp in $s2.Shapes) {
    if ($shp.HasTextFrame -eq -1 -and $shp.TextFrame.TextRange.Text -eq "Font End") {
        $shp.TextFrame.TextRange.Text = "Front End"
    }
}
